# Weekly refresh of the "Cebollín" (Hortaliza) price series at
# Vega Monumental Concepción: the oldest weekly record (row 129) is
# dropped, every remaining record shifts up one row, and a new record
# is appended at the bottom (row 142).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values for rows 129-142, columns D, I, J, K, L, M, N, O, P, Q
$rows = @(
    @{ Row = 129; D = 44943; I = "Primera"; J = 450; K = 2700; L = 2800; M = 2756; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 77;  Q = 36 },
    @{ Row = 130; D = 45216; I = "Primera"; J = 100; K = 4000; L = 4500; M = 4250; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 118; Q = 36 },
    @{ Row = 131; D = 45076; I = "Primera"; J = 100; K = 4000; L = 4200; M = 4100; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 114; Q = 36 },
    @{ Row = 132; D = 44595; I = "Primera"; J = 200; K = 600;  L = 700;  M = 650;  N = "`$/paquete 6 unidades";  O = "Región Metropolitana"; P = 108; Q = 6  },
    @{ Row = 133; D = 44692; I = "Primera"; J = 200; K = 600;  L = 700;  M = 650;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 108; Q = 6  },
    @{ Row = 134; D = 44692; I = "Segunda"; J = 100; K = 500;  L = 500;  M = 500;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 83;  Q = 6  },
    @{ Row = 135; D = 44901; I = "Primera"; J = 220; K = 3000; L = 3500; M = 3273; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 91;  Q = 36 },
    @{ Row = 136; D = 45092; I = "Primera"; J = 50;  K = 3500; L = 4000; M = 3800; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 106; Q = 36 },
    @{ Row = 137; D = 45190; I = "Primera"; J = 50;  K = 4500; L = 5000; M = 4700; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 131; Q = 36 },
    @{ Row = 138; D = 44965; I = "Primera"; J = 200; K = 700;  L = 800;  M = 750;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 125; Q = 6  },
    @{ Row = 139; D = 44965; I = "Segunda"; J = 100; K = 600;  L = 600;  M = 600;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 100; Q = 6  },
    @{ Row = 140; D = 44911; I = "Primera"; J = 200; K = 700;  L = 800;  M = 750;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 125; Q = 6  },
    @{ Row = 141; D = 44911; I = "Segunda"; J = 100; K = 600;  L = 600;  M = 600;  N = "`$/paquete 6 unidades";  O = "Región de Ñuble";      P = 100; Q = 6  },
    @{ Row = 142; D = 45267; I = "Primera"; J = 100; K = 4500; L = 4800; M = 4650; N = "`$/paquete 36 unidades"; O = "Región Metropolitana"; P = 129; Q = 36 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 4).Value  = $r.D   # D - Fecha
    $ws.Cells.Item($row, 9).Value  = $r.I   # I - Calidad
    $ws.Cells.Item($row, 10).Value = $r.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value = $r.K   # K - Precio mínimo
    $ws.Cells.Item($row, 12).Value = $r.L   # L - Precio máximo
    $ws.Cells.Item($row, 13).Value = $r.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 14).Value = $r.N   # N - Unidad de comercialización
    $ws.Cells.Item($row, 15).Value = $r.O   # O - Origen
    $ws.Cells.Item($row, 16).Value = $r.P   # P - Precio $/Kg
    $ws.Cells.Item($row, 17).Value = $r.Q   # Q - Kg o Unidades
}
